# Final edits to PIXm and minor update to MHDS
# -------------------------------------------------------------
# 1) Bump the fixed "datetimeFigureOut" footer date from
#    11/14/2019 -> 11/19/2019 on the slide master and every
#    slide layout (PowerPoint stores the cached text of that
#    field on each master/layout placeholder).
# 2) Update three "system" description textboxes that are
#    repeated (with slightly different wording) across the
#    MHDS partner-type diagrams, re-sizing the autofit textbox
#    that shrinks/grows as its text changes.
# -------------------------------------------------------------

$p = $ppt.ActivePresentation

function Update-DateField {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "11/14/2019") {
                $sh.TextFrame.TextRange.Text = "11/19/2019"
            }
        }
    }
}

# --- 1) Slide master + all custom (slide) layouts ----------------
$master = $p.Designs.Item(1).SlideMaster
Update-DateField $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateField $layouts.Item($i).Shapes
}

# --- 2) "Document Source" textbox (slide 4) -----------------------
$slide4 = $p.Slides.Item(4)
$src = $slide4.Shapes.Item("TextBox 26")
$src.TextFrame.TextRange.Paragraphs(1, 1).Text = "System that publishes Documents"

# --- 3) "Document Consumer" textbox (slide 5) ---------------------
# The wording changes AND one trailing blank paragraph is removed
# (the textbox keeps auto-fitting to its remaining lines).
$slide5 = $p.Slides.Item(5)
$cons = $slide5.Shapes.Item("TextBox 26")
$consRange = $cons.TextFrame.TextRange
$lastPara = $consRange.Paragraphs().Count
$consRange.Paragraphs($lastPara, 1).Delete()
$consRange.Paragraphs(1, 1).Text = "System that consumes documents"
$cons.Height = 356.2453

# --- 4) "Clinical Data Consumer" textbox (slide 6) ----------------
$slide6 = $p.Slides.Item(6)
$clin = $slide6.Shapes.Item("TextBox 26")
$clin.TextFrame.TextRange.Paragraphs(1, 1).Text = "System that consumes clinical data elements"
